$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 3) to the submissions table, mirroring the shape
# of row 2 but for a later submission (same person/camp/trip, later timestamp).

# A3 is an explicit empty "text" cell (matches A2's empty string), so we
# force text-typing via a leading apostrophe (Excel's normal "treat as text"
# idiom) rather than leaving it truly blank, which would omit the cell.
$ws.Range("A3").Value = "'"

$ws.Range("B3").Value = "احمد"

# C3 ("233") looks numeric, but the column stores it as text (same as C2),
# so use a leading apostrophe - Excel's normal "force text" idiom - to keep
# it a text value instead of letting it auto-coerce into a number.
$ws.Range("C3").Value = "'233"

$ws.Range("D3").Value = "الصمود"
$ws.Range("E3").Value = "الرحلة 2"
$ws.Range("F3").Value = "C2"
$ws.Range("G3").Value = "IDRF"
$ws.Range("H3").Value = "٠٥‏/٠٥‏/٢٠٢٥ ٠١:٥٩:١٤ م"
